$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row 1 with new columns P and Q (values 14 and 15), copying the format
# (bold, bordered, centered) from the existing header style used in O1.
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# For data rows 2-25: swap values in columns I/K and M/O, and add new columns P and Q with value 2
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I -> 2 (was 1)
    $ws.Cells.Item($r, 11).Value = 1  # K -> 1 (was 2)
    $ws.Cells.Item($r, 13).Value = 2  # M -> 2 (was 1)
    $ws.Cells.Item($r, 15).Value = 1  # O -> 1 (was 2)
    $ws.Cells.Item($r, 16).Value = 2  # P = 2
    $ws.Cells.Item($r, 17).Value = 2  # Q = 2
}
